$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 45970.01041666666
$ws.Cells.Item(2, 2).Value2 = 1322.946
$ws.Cells.Item(3, 1).Value2 = 45970.02083333334
$ws.Cells.Item(3, 2).Value2 = 1330.455
$ws.Cells.Item(4, 1).Value2 = 45970.03125
$ws.Cells.Item(4, 2).Value2 = 1337.381
$ws.Cells.Item(5, 1).Value2 = 45970.04166666666
$ws.Cells.Item(5, 2).Value2 = 1343.411
$ws.Cells.Item(6, 1).Value2 = 45970.05208333334
$ws.Cells.Item(6, 2).Value2 = 1362.475
$ws.Cells.Item(7, 1).Value2 = 45970.0625
$ws.Cells.Item(7, 2).Value2 = 1357.056
$ws.Cells.Item(8, 1).Value2 = 45970.07291666666
$ws.Cells.Item(8, 2).Value2 = 1349.069
$ws.Cells.Item(9, 1).Value2 = 45970.08333333334
$ws.Cells.Item(9, 2).Value2 = 1335.757
$ws.Cells.Item(10, 1).Value2 = 45970.09375
$ws.Cells.Item(10, 2).Value2 = 1312.52
$ws.Cells.Item(11, 1).Value2 = 45970.10416666666
$ws.Cells.Item(11, 2).Value2 = 1303.816
$ws.Cells.Item(12, 1).Value2 = 45970.11458333334
$ws.Cells.Item(12, 2).Value2 = 1288.597
$ws.Cells.Item(13, 1).Value2 = 45970.125
$ws.Cells.Item(13, 2).Value2 = 1267.889
$ws.Cells.Item(14, 1).Value2 = 45970.13541666666
$ws.Cells.Item(14, 2).Value2 = 1229.041
$ws.Cells.Item(15, 1).Value2 = 45970.14583333334
$ws.Cells.Item(15, 2).Value2 = 1208.372
$ws.Cells.Item(16, 1).Value2 = 45970.15625
$ws.Cells.Item(16, 2).Value2 = 1201.594
$ws.Cells.Item(17, 1).Value2 = 45970.16666666666
$ws.Cells.Item(17, 2).Value2 = 1182.459
$ws.Cells.Item(18, 1).Value2 = 45970.17708333334
$ws.Cells.Item(18, 2).Value2 = 1129.25
$ws.Cells.Item(19, 1).Value2 = 45970.1875
$ws.Cells.Item(19, 2).Value2 = 1098.795
$ws.Cells.Item(20, 1).Value2 = 45970.19791666666
$ws.Cells.Item(20, 2).Value2 = 1079.545
$ws.Cells.Item(21, 1).Value2 = 45970.20833333334
$ws.Cells.Item(21, 2).Value2 = 1060.094
$ws.Cells.Item(22, 1).Value2 = 45970.21875
$ws.Cells.Item(22, 2).Value2 = 1001.651
$ws.Cells.Item(23, 1).Value2 = 45970.22916666666
$ws.Cells.Item(23, 2).Value2 = 983.398
$ws.Cells.Item(24, 1).Value2 = 45970.23958333334
$ws.Cells.Item(24, 2).Value2 = 965.845
$ws.Cells.Item(25, 1).Value2 = 45970.25
$ws.Cells.Item(25, 2).Value2 = 951.062
$ws.Cells.Item(26, 1).Value2 = 45970.26041666666
$ws.Cells.Item(26, 2).Value2 = 860.994
$ws.Cells.Item(27, 1).Value2 = 45970.27083333334
$ws.Cells.Item(27, 2).Value2 = 847.509
$ws.Cells.Item(28, 1).Value2 = 45970.28125
$ws.Cells.Item(28, 2).Value2 = 826.183
$ws.Cells.Item(29, 1).Value2 = 45970.29166666666
$ws.Cells.Item(29, 2).Value2 = 812.2809999999999
$ws.Cells.Item(30, 1).Value2 = 45970.30208333334
$ws.Cells.Item(30, 2).Value2 = 713.424
$ws.Cells.Item(31, 1).Value2 = 45970.3125
$ws.Cells.Item(31, 2).Value2 = 689.9349999999999
$ws.Cells.Item(32, 1).Value2 = 45970.32291666666
$ws.Cells.Item(32, 2).Value2 = 685.6799999999999
$ws.Cells.Item(33, 1).Value2 = 45970.33333333334
$ws.Cells.Item(33, 2).Value2 = 660.442
$ws.Cells.Item(34, 1).Value2 = 45970.34375
$ws.Cells.Item(34, 2).Value2 = 575.135
$ws.Cells.Item(35, 1).Value2 = 45970.35416666666
$ws.Cells.Item(35, 2).Value2 = 563.577
$ws.Cells.Item(36, 1).Value2 = 45970.36458333334
$ws.Cells.Item(36, 2).Value2 = 543.49
$ws.Cells.Item(37, 1).Value2 = 45970.375
$ws.Cells.Item(37, 2).Value2 = 525.671
$ws.Cells.Item(38, 1).Value2 = 45970.38541666666
$ws.Cells.Item(38, 2).Value2 = 439.411
$ws.Cells.Item(39, 1).Value2 = 45970.39583333334
$ws.Cells.Item(39, 2).Value2 = 421.232
$ws.Cells.Item(40, 1).Value2 = 45970.40625
$ws.Cells.Item(40, 2).Value2 = 408.565
$ws.Cells.Item(41, 1).Value2 = 45970.41666666666
$ws.Cells.Item(41, 2).Value2 = 396.959
$ws.Cells.Item(42, 1).Value2 = 45970.42708333334
$ws.Cells.Item(42, 2).Value2 = 338.414
$ws.Cells.Item(43, 1).Value2 = 45970.4375
$ws.Cells.Item(43, 2).Value2 = 333.146
$ws.Cells.Item(44, 1).Value2 = 45970.44791666666
$ws.Cells.Item(44, 2).Value2 = 327.795
$ws.Cells.Item(45, 1).Value2 = 45970.45833333334
$ws.Cells.Item(45, 2).Value2 = 320.928
$ws.Cells.Item(46, 1).Value2 = 45970.46875
$ws.Cells.Item(46, 2).Value2 = 298.231
$ws.Cells.Item(47, 1).Value2 = 45970.47916666666
$ws.Cells.Item(47, 2).Value2 = 294.683
$ws.Cells.Item(48, 1).Value2 = 45970.48958333334
$ws.Cells.Item(48, 2).Value2 = 290.608
$ws.Cells.Item(49, 1).Value2 = 45970.5
$ws.Cells.Item(49, 2).Value2 = 286.349
$ws.Cells.Item(50, 1).Value2 = 45970.51041666666
$ws.Cells.Item(50, 2).Value2 = 254.591
$ws.Cells.Item(51, 1).Value2 = 45970.52083333334
$ws.Cells.Item(51, 2).Value2 = 250.447
$ws.Cells.Item(52, 1).Value2 = 45970.53125
$ws.Cells.Item(52, 2).Value2 = 243.81
$ws.Cells.Item(53, 1).Value2 = 45970.54166666666
$ws.Cells.Item(53, 2).Value2 = 239.759
$ws.Cells.Item(54, 1).Value2 = 45970.55208333334
$ws.Cells.Item(54, 2).Value2 = 208.813
$ws.Cells.Item(55, 1).Value2 = 45970.5625
$ws.Cells.Item(55, 2).Value2 = 204.906
$ws.Cells.Item(56, 1).Value2 = 45970.57291666666
$ws.Cells.Item(56, 2).Value2 = 201.216
$ws.Cells.Item(57, 1).Value2 = 45970.58333333334
$ws.Cells.Item(57, 2).Value2 = 197.736
$ws.Cells.Item(58, 1).Value2 = 45970.59375
$ws.Cells.Item(58, 2).Value2 = 182.476
$ws.Cells.Item(59, 1).Value2 = 45970.60416666666
$ws.Cells.Item(59, 2).Value2 = 179.577
$ws.Cells.Item(60, 1).Value2 = 45970.61458333334
$ws.Cells.Item(60, 2).Value2 = 177.85
$ws.Cells.Item(61, 1).Value2 = 45970.625
$ws.Cells.Item(61, 2).Value2 = 176.391
$ws.Cells.Item(62, 1).Value2 = 45970.63541666666
$ws.Cells.Item(62, 2).Value2 = 171.384
$ws.Cells.Item(63, 1).Value2 = 45970.64583333334
$ws.Cells.Item(63, 2).Value2 = 174.06
$ws.Cells.Item(64, 1).Value2 = 45970.65625
$ws.Cells.Item(64, 2).Value2 = 176.973
$ws.Cells.Item(65, 1).Value2 = 45970.66666666666
$ws.Cells.Item(65, 2).Value2 = 179.983
$ws.Cells.Item(66, 1).Value2 = 45970.67708333334
$ws.Cells.Item(66, 2).Value2 = 183.646
$ws.Cells.Item(67, 1).Value2 = 45970.6875
$ws.Cells.Item(67, 2).Value2 = 187.258
$ws.Cells.Item(68, 1).Value2 = 45970.69791666666
$ws.Cells.Item(68, 2).Value2 = 190.854
$ws.Cells.Item(69, 1).Value2 = 45970.70833333334
$ws.Cells.Item(69, 2).Value2 = 194.335
$ws.Cells.Item(70, 1).Value2 = 45970.71875
$ws.Cells.Item(70, 2).Value2 = 196.092
$ws.Cells.Item(71, 1).Value2 = 45970.72916666666
$ws.Cells.Item(71, 2).Value2 = 196.927
$ws.Cells.Item(72, 1).Value2 = 45970.73958333334
$ws.Cells.Item(72, 2).Value2 = 197.585
$ws.Cells.Item(73, 1).Value2 = 45970.75
$ws.Cells.Item(73, 2).Value2 = 198.126
$ws.Cells.Item(74, 1).Value2 = 45970.76041666666
$ws.Cells.Item(74, 2).Value2 = 201.43
$ws.Cells.Item(75, 1).Value2 = 45970.77083333334
$ws.Cells.Item(75, 2).Value2 = 201.684
$ws.Cells.Item(76, 1).Value2 = 45970.78125
$ws.Cells.Item(76, 2).Value2 = 201.132
$ws.Cells.Item(77, 1).Value2 = 45970.79166666666
$ws.Cells.Item(77, 2).Value2 = 200.322
$ws.Cells.Item(78, 1).Value2 = 45970.80208333334
$ws.Cells.Item(78, 2).Value2 = 191.729
$ws.Cells.Item(79, 1).Value2 = 45970.8125
$ws.Cells.Item(79, 2).Value2 = 189.691
$ws.Cells.Item(80, 1).Value2 = 45970.82291666666
$ws.Cells.Item(80, 2).Value2 = 187.817
$ws.Cells.Item(81, 1).Value2 = 45970.83333333334
$ws.Cells.Item(81, 2).Value2 = 185.464
$ws.Cells.Item(82, 1).Value2 = 45970.84375
$ws.Cells.Item(82, 2).Value2 = 171.915
$ws.Cells.Item(83, 1).Value2 = 45970.85416666666
$ws.Cells.Item(83, 2).Value2 = 169.076
$ws.Cells.Item(84, 1).Value2 = 45970.86458333334
$ws.Cells.Item(84, 2).Value2 = 165.77
$ws.Cells.Item(85, 1).Value2 = 45970.875
$ws.Cells.Item(85, 2).Value2 = 162.154
$ws.Cells.Item(86, 1).Value2 = 45970.88541666666
$ws.Cells.Item(86, 2).Value2 = 155.216
$ws.Cells.Item(87, 1).Value2 = 45970.89583333334
$ws.Cells.Item(87, 2).Value2 = 151.082
$ws.Cells.Item(88, 1).Value2 = 45970.90625
$ws.Cells.Item(88, 2).Value2 = 147.821
$ws.Cells.Item(89, 1).Value2 = 45970.91666666666
$ws.Cells.Item(89, 2).Value2 = 144.626
$ws.Cells.Item(90, 1).Value2 = 45970.92708333334
$ws.Cells.Item(90, 2).Value2 = 128.292
$ws.Cells.Item(91, 1).Value2 = 45970.9375
$ws.Cells.Item(91, 2).Value2 = 124.73
$ws.Cells.Item(92, 1).Value2 = 45970.94791666666
$ws.Cells.Item(92, 2).Value2 = 121.594
$ws.Cells.Item(93, 1).Value2 = 45970.95833333334
$ws.Cells.Item(93, 2).Value2 = 119.13
$ws.Cells.Item(94, 1).Value2 = 45970.96875
$ws.Cells.Item(94, 2).Value2 = 0
$ws.Cells.Item(95, 1).Value2 = 45970.97916666666
$ws.Cells.Item(95, 2).Value2 = 0
$ws.Cells.Item(96, 1).Value2 = 45970.98958333334
$ws.Cells.Item(96, 2).Value2 = 0
$ws.Cells.Item(97, 1).Value2 = 45971
$ws.Cells.Item(97, 2).Value2 = 0
